$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 194 - 四方坪站充电量(kw) - 2025-12-06 (Excel serial 45997)
$ws.Range("A194").Value = 45997
$ws.Range("B194").Value = "四方坪站充电量(kw)"

$row194 = @(637.74700000000007, 1449.9689999999998, 548.71299999999985, 625.59700000000009, 282.44499999999999, 630.20500000000004, 495.15699999999993, 230.077, 182.47800000000001, 133.375, 135.42000000000002, 358.44299999999998, 734.65, 1526.9880000000001, 668.6930000000001, 580.06999999999994, 438.36600000000004, 153.953, 131.45099999999999, 13.698999999999998, 110.83399999999999, 22.159999999999997, 56.410000000000004, 9.48)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "194").Value = $row194[$i]
}

# Row 195 - 高岭站充电量(kw) - 2025-12-06 (Excel serial 45997)
$ws.Range("A195").Value = 45997
$ws.Range("B195").Value = "高岭站充电量(kw)"

$row195 = @(334.43499999999995, 582.17199999999991, 232.34899999999999, 34.950000000000003, 137.54499999999999, 224.708, 216.04499999999996, 121.6, 141.494, 146.69200000000001, 214.57799999999997, 342.71000000000004, 252.601, 311.56799999999998, 295.90800000000002, 418.28500000000003, 244.404, 172.61300000000003, 95.135999999999996, 37.085999999999999, 0, 61.686999999999998, 77.411000000000001, 60.108999999999995)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "195").Value = $row195[$i]
}

# Update the view: selection moves to E199, top-left stays at A169
$ws.Range("E199").Select() | Out-Null
